$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista Atividades")

# Remove the stored cell selection on the sheet view (B6) by re-selecting A1
$ws.Activate()
$ws.Range("A1").Select()

# New rows of data appended to the activity list
$ws.Range("A24").Value = "Cadastro de Docente"
$ws.Range("B24").Value = "Ademar Júnior/Jean Lucas/Fernando Gonzaga"
$ws.Range("C24").Value = "SIM"

$ws.Range("A25").Value = "Cadastro de Atividades"
$ws.Range("B25").Value = "Mário Hayasaki/Guilherme Moreno"
$ws.Range("C25").Value = "SIM"

# Apply the same bordered style as the rest of the table (reuse existing
# border style from row 23), plus an explicit (but empty) fill application
$ws.Range("A23:C23").Copy()
$ws.Range("A24:C25").PasteSpecial(-4122)
$ws.Range("A24:C25").Interior.Pattern = -4142
